# Applies the "Visão Emergente" edit:
#  1. Split the opening paragraph "Ideias: Visão Emergente." into two
#     paragraphs: "Visão Emergente" and the new explanatory sentence.
#  2. Split the run of the "Lâmpada Inteligente..." bullet after "Lâmpa"
#     and move the _GoBack bookmark there.
#  3. Split the run of the "Micro Modem 3D..." bullet into three runs,
#     turning the "D" into "G" along the way.
#  4. The stray _GoBack bookmark that used to sit at the very end of the
#     document (after the "Geladeira..." bullet) disappears automatically
#     because _GoBack is unique and gets relocated in step 2.

$d = $word.ActiveDocument

# --- 1. Title paragraph split -------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$r1.Text = "Visão Emergente"
$r1.InsertParagraphAfter()

$p2 = $d.Paragraphs.Item(2)
$p2.Range.Text = "No processo visionário a visão emergente é o ato de buscar ideias para a criação do produto ou serviço. Neste momento surgiram as seguintes ideias:"

# --- 2. "Lâmpada Inteligente..." bullet: split after "Lâmpa" and plant --
# --- the _GoBack bookmark at the split point (removing it from wherever
# --- it used to be, since a document can only have one _GoBack). -------
$pLampada = $d.Paragraphs.Item(3)
$lampadaStart = $pLampada.Range.Start
$splitPoint = $lampadaStart + 5   # "Lâmpa" is 5 characters
$splitRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $splitRange) | Out-Null

# --- 3. "Micro Modem 3D..." bullet: "D" -> "G" and split into 3 runs ---
$pModem = $d.Paragraphs.Item(7)
$modemStart = $pModem.Range.Start

# "Micro Modem 3" is 13 characters; the following character is the "D".
$dRange = $d.Range($modemStart + 13, $modemStart + 14)
$dRange.Text = "G"

# Force the run boundaries by momentarily bookmarking (and immediately
# removing) the two split points; the text stays split into separate
# runs even after the temporary bookmark is deleted.
$split1 = $d.Range($modemStart + 13, $modemStart + 13)
$tmp1 = $d.Bookmarks.Add("TempSplitA", $split1)
$d.Bookmarks.Item("TempSplitA").Delete()

$split2 = $d.Range($modemStart + 14, $modemStart + 14)
$tmp2 = $d.Bookmarks.Add("TempSplitB", $split2)
$d.Bookmarks.Item("TempSplitB").Delete()
